# Updates cryptos list values (Price in column D, Volume(1h) % in column E)
# on the active worksheet, matching refreshed source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.612.21"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "'1.926.81"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'326.60"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D7").Value = "'0.4828"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'0.4065"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'0.08241"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D11").Value = "'23.88"
$ws.Range("D12").Value = "'1.929.96"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "'6.116"
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").Value = "'7.281"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").Value = "'0.06865"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").Value = "'29.606.23"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "'5.696"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "'12.02"
$ws.Range("E23").Value = "  +2.29%  "
$ws.Range("D24").Value = "'2.189"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").Value = "'2.166.70"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").Value = "'156.30"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "'6.443"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").Value = "'2.098"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "'120.88"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").Value = "'1.017"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "'0.09650"
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("D33").Value = "'5.633"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("D34").Value = "'3.576"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").Value = "'1.381"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").Value = "'0.06395"
$ws.Range("E36").Value = "  +5.00%  "
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("D38").Value = "'1.188"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("D39").Value = "'0.5968"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").Value = "'7.906"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "'0.1856"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").Value = "'2.433"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "'12.47"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "'0.07551"
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("D47").Value = "'0.5577"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "'1.972"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").Value = "'119.66"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").Value = "'2.440"
$ws.Range("E50").Value = "  +3.62%  "
$ws.Range("D51").Value = "'72.25"
$ws.Range("E51").Value = "  -0.30%  "
